# "break out stock.yaml completed" -- chartink_screener.xlsx, "10per change" sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("10per change")

# --- E58 / E59: bsecode had been scraped as text; re-store as real numbers. ---
$ws.Range("E58").Value = 590024
$ws.Range("E59").Value = 543220

# --- Append the four newly-scraped rows (60-63). ---
$newRows = @(
    @{ Row = 60; DateTime = "26/06/2024 05:45:37"; Sr = 1; NseCode = "SUPREMEIND"; Name = "Supreme Industries Limited";                   BseCode = "509930"; PerChg = -1.41; Close = 5804.65; Volume = 25306 },
    @{ Row = 61; DateTime = "26/06/2024 05:45:37"; Sr = 2; NseCode = "FACT";       Name = "Fertilizers And Chemicals Travancore Limited"; BseCode = "590024"; PerChg = 2.25;  Close = 1021.55; Volume = 968379 },
    @{ Row = 62; DateTime = "26/06/2024 05:45:37"; Sr = 3; NseCode = "MAXHEALTH";  Name = "Max Healthcare Institute Ltd";                  BseCode = "543220"; PerChg = -2.06; Close = 875.05;  Volume = 905629 },
    @{ Row = 63; DateTime = "26/06/2024 05:45:37"; Sr = 4; NseCode = "GAIL";       Name = "Gail (india) Limited";                          BseCode = "532155"; PerChg = -1.57; Close = 209.8;   Volume = 4986674 }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = $r.DateTime
    $ws.Cells.Item($row, 2).Value = $r.Sr
    $ws.Cells.Item($row, 3).Value = $r.NseCode
    $ws.Cells.Item($row, 4).Value = $r.Name

    # bsecode stays TEXT for these rows (matches the scraper's raw output) --
    # go through a TEXT() formula and paste the result back as a value so the
    # numeric-looking string is NOT coerced into a number, without leaving a
    # "@" number-format behind on the cell.
    $cell = $ws.Cells.Item($row, 5)
    $cell.Formula = '=TEXT(' + $r.BseCode + ',"0")'
    $cell.Copy()
    $cell.PasteSpecial(-4163)

    $ws.Cells.Item($row, 6).Value = $r.PerChg
    $ws.Cells.Item($row, 7).Value = $r.Close
    $ws.Cells.Item($row, 8).Value = $r.Volume
}
